$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Rows.Item(59).Insert()
$ws.Range("A59:D59").Clear()

$ws.Range("H60").Copy()
$ws.Range("H59").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I60").Copy()
$ws.Range("I59").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("K60").Copy()
$ws.Range("K59").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A192").Value = "x"
$ws.Range("B192").Value = "y"
$ws.Range("C192").Value = 15000
$ws.Range("D192").Formula = "=D191+C192"

Write-Host "D192 formula: $($ws.Range("D192").Formula)"
Write-Host "D192 value: $($ws.Range("D192").Value2)"
